# Add November (through 11-01) 2021 data and roll the running totals forward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook/sheet name tracks the "through" date of the data.
$ws.Name = "Through 2021-11-01"

# October's row no longer needs the "(through 10-31)" qualifier now that the
# month is complete.
$ws.Range("A12").Value = "October"

# Insert a new row for November's (partial) data; this pushes the old
# "Total" row from 13 down to 14.
$ws.Rows.Item(13).Insert()

# Give the new row's label cell (A13) the same look (bold, bordered) as the
# other month-label cells by copying the format from the row above.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the November (through 11-01) row.
$ws.Range("A13").Value = "November (through 11-01)"
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 3
$ws.Range("I13").Value = 5
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 3
$ws.Range("M13").Value = 0.25
$ws.Range("O13").Value = 1
$ws.Range("R13").Value = 9
$ws.Range("U13").Value = 5

# The row-insert auto-filled the percent-format columns (which carry a
# column-level style) with empty formatted cells; November has no data for
# those year/arrest-rate columns, so drop those placeholder cells.
$ws.Range("D13").Clear()
$ws.Range("G13").Clear()
$ws.Range("J13").Clear()
$ws.Range("P13").Clear()
$ws.Range("S13").Clear()
$ws.Range("V13").Clear()

# Update the running totals, now on row 14, to include November's data.
$ws.Range("C14").Value = 227
$ws.Range("D14").Value = 0.1236
$ws.Range("F14").Value = 437
$ws.Range("G14").Value = 0.1063
$ws.Range("I14").Value = 654
$ws.Range("J14").Value = 0.0853
$ws.Range("K14").Value = 68
$ws.Range("L14").Value = 551
$ws.Range("M14").Value = 0.1099
$ws.Range("O14").Value = 435
$ws.Range("P14").Value = 0.0994
$ws.Range("R14").Value = 1012
$ws.Range("S14").Value = 0.0507
$ws.Range("U14").Value = 1364
$ws.Range("V14").Value = 0.0587

# Column A is slightly wider to fit the longer "November (through 11-01)"
# label (target raw width 24.7109375; the COM layer here quantizes
# ColumnWidth to whole pixels, so 23.8 is the closest achievable input).
$ws.Columns.Item(1).ColumnWidth = 23.8

Write-Output "Applied November 2021-11-01 update"
